# Refresh the crypto price / volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe is the normal Excel "treat as text" quote-prefix, so
# numeric-looking Price values (e.g. "1.002", "0.9995") stay text cells, same
# as the rest of the column, instead of being auto-parsed into numbers.

$ws.Range("D2").Value = '23.436.63'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.643.76'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '''0.9995'
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("D6").Value = '''299.77'
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("D7").Value = '''0.3795'
$ws.Range("E7").Value = '  -1.20%  '

$ws.Range("D8").Value = '''50.54'
$ws.Range("E8").Value = '  -0.94%  '

$ws.Range("D9").Value = '''0.3505'
$ws.Range("E9").Value = '  -2.49%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.08060'
$ws.Range("E10").Value = '  -1.33%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '''1.217'
$ws.Range("E11").Value = '  -0.95%  '

$ws.Range("D12").Value = '''1.002'
$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("D13").Value = '''22.04'
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").Value = '''6.307'
$ws.Range("E14").Value = '  -2.13%  '

$ws.Range("D15").Value = '''7.250'
$ws.Range("E15").Value = '  -2.80%  '

$ws.Range("D16").Value = '''0.00001209'
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").Value = '1.639.27'
$ws.Range("E17").Value = '  -0.91%  '

$ws.Range("D18").Value = '''94.94'
$ws.Range("E18").Value = '  -2.59%  '

$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").Value = '''6.618'
$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").Value = '''17.41'
$ws.Range("E21").Value = '  -0.72%  '

$ws.Range("D22").Value = '''0.9991'
$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("D23").Value = '''12.44'
$ws.Range("E23").Value = '  -1.28%  '

$ws.Range("D24").Value = '23.459.37'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '''2.416'
$ws.Range("E25").Value = '  -2.65%  '

$ws.Range("D26").Value = '''2.963'
$ws.Range("E26").Value = '  -2.33%  '

$ws.Range("D27").Value = '''21.05'
$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("D28").Value = '''151.00'
$ws.Range("E28").Value = '  -1.73%  '

$ws.Range("D29").Value = '''5.171'
$ws.Range("E29").Value = '  -1.05%  '

$ws.Range("D30").Value = '''132.12'
$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("D31").Value = '1.832.02'
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("D32").Value = '''6.826'
$ws.Range("E32").Value = '  -4.06%  '

$ws.Range("D33").Value = '''2.137'
$ws.Range("E33").Value = '  -4.83%  '

$ws.Range("D34").Value = '''11.17'
$ws.Range("E34").Value = '  -8.38%  '

$ws.Range("D35").Value = '''0.9861'
$ws.Range("E35").Value = '  -6.45%  '

$ws.Range("D36").Value = '''0.02693'
$ws.Range("E36").Value = '  -3.43%  '

$ws.Range("D37").Value = '''0.08784'
$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("D38").Value = '''5.908'
$ws.Range("E38").Value = '  -2.62%  '

$ws.Range("D39").Value = '''0.2418'
$ws.Range("E39").Value = '  -3.28%  '

$ws.Range("D40").Value = '''0.06767'
$ws.Range("E40").Value = '  -3.02%  '

$ws.Range("D41").Value = '''12.79'
$ws.Range("E41").Value = '  -2.01%  '

$ws.Range("D42").Value = '''0.6853'
$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("D43").Value = '''1.294'
$ws.Range("E43").Value = '  -2.88%  '

$ws.Range("D44").Value = '''15.50'
$ws.Range("E44").Value = '  -2.40%  '

$ws.Range("D45").Value = '''0.9984'
$ws.Range("E45").Value = '  -0.43%  '

$ws.Range("D46").Value = '''0.6366'
$ws.Range("E46").Value = '  -2.05%  '

$ws.Range("D47").Value = '''2.245'
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").Value = '''3.918'
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''127.09'
$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.07675'
$ws.Range("E50").Value = '  -2.47%  '

$ws.Range("D51").Value = '''1.231'
$ws.Range("E51").Value = '  +2.40%  '
